$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 116789.285257321
$ws.Cells.Item(2, 3).Value = 109504.654090983
$ws.Cells.Item(2, 4).Value = 105967.282430223
$ws.Cells.Item(2, 5).Value = 124686.275337719
$ws.Cells.Item(2, 6).Value = 127602.299021983
$ws.Cells.Item(2, 9).Value = 25763.2852573207
$ws.Cells.Item(3, 2).Value = 103994.843739405
$ws.Cells.Item(3, 3).Value = 96970.4852970524
$ws.Cells.Item(3, 4).Value = 91498.9115345899
$ws.Cells.Item(3, 5).Value = 111911.78536891
$ws.Cells.Item(3, 6).Value = 114452.468094831
$ws.Cells.Item(3, 9).Value = 52488.8437394046
$ws.Cells.Item(4, 2).Value = 122358.327180849
$ws.Cells.Item(4, 3).Value = 114390.364234689
$ws.Cells.Item(4, 4).Value = 109100.579549827
$ws.Cells.Item(4, 5).Value = 130316.017589693
$ws.Cells.Item(4, 6).Value = 134395.83284113
$ws.Cells.Item(4, 9).Value = 34208.3271808492
$ws.Cells.Item(5, 2).Value = 113378.860349499
$ws.Cells.Item(5, 3).Value = 104858.729261052
$ws.Cells.Item(5, 4).Value = 100967.506649076
$ws.Cells.Item(5, 5).Value = 121801.160779743
$ws.Cells.Item(5, 6).Value = 125233.73928071
$ws.Cells.Item(5, 9).Value = 12116.8603494988
$ws.Cells.Item(6, 2).Value = 114300.964768311
$ws.Cells.Item(6, 3).Value = 105382.142403267
$ws.Cells.Item(6, 4).Value = 101167.028048125
$ws.Cells.Item(6, 5).Value = 123469.934343354
$ws.Cells.Item(6, 6).Value = 127359.072437067
$ws.Cells.Item(6, 9).Value = 16649.9647683114
$ws.Cells.Item(7, 2).Value = 108838.650428902
$ws.Cells.Item(7, 3).Value = 99226.3535260441
$ws.Cells.Item(7, 4).Value = 94465.7811819188
$ws.Cells.Item(7, 5).Value = 118313.481149257
$ws.Cells.Item(7, 6).Value = 122579.903555989
$ws.Cells.Item(7, 9).Value = 9519.65042890183
$ws.Cells.Item(8, 2).Value = 112901.749864161
$ws.Cells.Item(8, 3).Value = 103080.742812283
$ws.Cells.Item(8, 4).Value = 98608.2903223724
$ws.Cells.Item(8, 5).Value = 123094.560429769
$ws.Cells.Item(8, 6).Value = 127485.546552157
$ws.Cells.Item(8, 9).Value = 10597.7498641613
$ws.Cells.Item(9, 2).Value = 111734.638687349
$ws.Cells.Item(9, 3).Value = 101793.945128848
$ws.Cells.Item(9, 4).Value = 97674.4120607763
$ws.Cells.Item(9, 5).Value = 121776.415535447
$ws.Cells.Item(9, 6).Value = 126736.65341991
$ws.Cells.Item(9, 9).Value = 5599.63868734882
$ws.Cells.Item(10, 2).Value = 102440.98453939
$ws.Cells.Item(10, 3).Value = 91169.5008043546
$ws.Cells.Item(10, 4).Value = 86228.2807855309
$ws.Cells.Item(10, 5).Value = 112844.595725229
$ws.Cells.Item(10, 6).Value = 118279.212215841
$ws.Cells.Item(10, 9).Value = -2936.01546061043
$ws.Cells.Item(11, 2).Value = 105296.866448337
$ws.Cells.Item(11, 3).Value = 94654.0507735155
$ws.Cells.Item(11, 4).Value = 89645.6189910836
$ws.Cells.Item(11, 5).Value = 116066.184623402
$ws.Cells.Item(11, 6).Value = 120397.303695727
$ws.Cells.Item(11, 9).Value = 9663.86644833727
$ws.Cells.Item(12, 2).Value = 104183.38408271
$ws.Cells.Item(12, 3).Value = 92437.722982702
$ws.Cells.Item(12, 4).Value = 86585.137141834
$ws.Cells.Item(12, 5).Value = 115962.989387261
$ws.Cells.Item(12, 6).Value = 121294.564225157
$ws.Cells.Item(12, 9).Value = 3622.38408271036
$ws.Cells.Item(13, 2).Value = 101956.498228134
$ws.Cells.Item(13, 3).Value = 90688.8159788464
$ws.Cells.Item(13, 4).Value = 86262.7232827484
$ws.Cells.Item(13, 5).Value = 113591.101420175
$ws.Cells.Item(13, 6).Value = 117998.81434705
$ws.Cells.Item(13, 9).Value = 1747.49822813405
$ws.Cells.Item(14, 2).Value = 119363.067367022
$ws.Cells.Item(14, 3).Value = 107113.87418771
$ws.Cells.Item(14, 4).Value = 100328.898082453
$ws.Cells.Item(14, 5).Value = 132636.75124848
$ws.Cells.Item(14, 6).Value = 138926.083985212
$ws.Cells.Item(14, 9).Value = 21106.067367022
$ws.Cells.Item(15, 2).Value = 106592.043663742
$ws.Cells.Item(15, 3).Value = 94093.0347442018
$ws.Cells.Item(15, 4).Value = 86744.864128847
$ws.Cells.Item(15, 5).Value = 119690.525147055
$ws.Cells.Item(15, 6).Value = 124407.57938035
$ws.Cells.Item(15, 9).Value = 23578.0436637417
$ws.Cells.Item(16, 2).Value = 124766.571323058
$ws.Cells.Item(16, 3).Value = 113511.217972994
$ws.Cells.Item(16, 4).Value = 104946.335496161
$ws.Cells.Item(16, 5).Value = 139381.549180683
$ws.Cells.Item(16, 6).Value = 144475.482368075
$ws.Cells.Item(16, 9).Value = 5343.57132305781
$ws.Cells.Item(17, 2).Value = 115612.264832069
$ws.Cells.Item(17, 3).Value = 104214.689117265
$ws.Cells.Item(17, 4).Value = 95996.5849916201
$ws.Cells.Item(17, 5).Value = 128076.387858619
$ws.Cells.Item(17, 6).Value = 135557.186356835
$ws.Cells.Item(17, 9).Value = 5227.2648320693
$ws.Cells.Item(18, 2).Value = 116600.750621209
$ws.Cells.Item(18, 3).Value = 103392.792946294
$ws.Cells.Item(18, 4).Value = 95489.8857375593
$ws.Cells.Item(18, 5).Value = 130333.418268387
$ws.Cells.Item(18, 6).Value = 138049.195714955
$ws.Cells.Item(18, 9).Value = 11207.7506212087
$ws.Cells.Item(19, 2).Value = 111081.779036424
$ws.Cells.Item(19, 3).Value = 96423.5152660892
$ws.Cells.Item(19, 4).Value = 87936.4745368495
$ws.Cells.Item(19, 5).Value = 125321.46513576
$ws.Cells.Item(19, 6).Value = 133704.547811
$ws.Cells.Item(19, 9).Value = 9448.77903642383
$ws.Cells.Item(20, 2).Value = 115361.184541697
$ws.Cells.Item(20, 3).Value = 100674.519507487
$ws.Cells.Item(20, 4).Value = 95060.5783766874
$ws.Cells.Item(20, 5).Value = 130145.933764809
$ws.Cells.Item(20, 6).Value = 137590.760249334
$ws.Cells.Item(20, 9).Value = 3074.1845416974
$ws.Cells.Item(21, 2).Value = 114740.647141269
$ws.Cells.Item(21, 3).Value = 100016.830064625
$ws.Cells.Item(21, 4).Value = 93392.102368367
$ws.Cells.Item(21, 5).Value = 130997.699663737
$ws.Cells.Item(21, 6).Value = 137636.574409553
$ws.Cells.Item(21, 9).Value = 12679.6471412692
$ws.Cells.Item(22, 2).Value = 104952.615386672
$ws.Cells.Item(22, 3).Value = 89674.9414246369
$ws.Cells.Item(22, 4).Value = 82947.0140564044
$ws.Cells.Item(22, 5).Value = 121938.654705538
$ws.Cells.Item(22, 6).Value = 129746.938891176
$ws.Cells.Item(22, 9).Value = 3251.61538667229
$ws.Cells.Item(23, 2).Value = 107561.807086219
$ws.Cells.Item(23, 3).Value = 92351.5068210012
$ws.Cells.Item(23, 4).Value = 84638.1722618577
$ws.Cells.Item(23, 5).Value = 123146.25536691
$ws.Cells.Item(23, 6).Value = 133689.200731027
$ws.Cells.Item(23, 9).Value = 11619.807086219
$ws.Cells.Item(24, 2).Value = 106712.988694516
$ws.Cells.Item(24, 3).Value = 91456.103131387
$ws.Cells.Item(24, 4).Value = 81933.9759173043
$ws.Cells.Item(24, 5).Value = 123012.004939941
$ws.Cells.Item(24, 6).Value = 133746.564409717
$ws.Cells.Item(24, 9).Value = 4114.98869451611
$ws.Cells.Item(25, 2).Value = 104385.70881101
$ws.Cells.Item(25, 3).Value = 89043.7445456703
$ws.Cells.Item(25, 4).Value = 79602.8833529934
$ws.Cells.Item(25, 5).Value = 120976.675233669
$ws.Cells.Item(25, 6).Value = 131869.358714816
$ws.Cells.Item(25, 9).Value = -3771.29118899028
$ws.Cells.Item(26, 2).Value = 122141.477911815
$ws.Cells.Item(26, 3).Value = 105534.666587195
$ws.Cells.Item(26, 4).Value = 95374.3861371691
$ws.Cells.Item(26, 5).Value = 139558.40178172
$ws.Cells.Item(26, 6).Value = 152981.36634396
$ws.Cells.Item(26, 9).Value = 12610.4779118146
$ws.Cells.Item(27, 2).Value = 109561.797553954
$ws.Cells.Item(27, 3).Value = 92780.5522298092
$ws.Cells.Item(27, 4).Value = 85886.3398778998
$ws.Cells.Item(27, 5).Value = 126090.426765374
$ws.Cells.Item(27, 6).Value = 137712.997675596
$ws.Cells.Item(27, 9).Value = 12774.7975539543
$ws.Cells.Item(28, 2).Value = 127252.841874435
$ws.Cells.Item(28, 3).Value = 110887.877433386
$ws.Cells.Item(28, 4).Value = 101395.126769917
$ws.Cells.Item(28, 5).Value = 145766.396634774
$ws.Cells.Item(28, 6).Value = 156413.742412195
$ws.Cells.Item(28, 9).Value = 13249.8418744346
$ws.Cells.Item(29, 2).Value = 118715.747077947
$ws.Cells.Item(29, 3).Value = 100236.909221468
$ws.Cells.Item(29, 4).Value = 91640.9319959471
$ws.Cells.Item(29, 5).Value = 137037.959709605
$ws.Cells.Item(29, 6).Value = 147394.597461622
$ws.Cells.Item(29, 9).Value = 22782.7470779472
$ws.Cells.Item(30, 2).Value = 119910.734740832
$ws.Cells.Item(30, 3).Value = 103679.236224741
$ws.Cells.Item(30, 4).Value = 94063.7443721575
$ws.Cells.Item(30, 5).Value = 137524.877157735
$ws.Cells.Item(30, 6).Value = 151094.142225869
$ws.Cells.Item(30, 9).Value = 16998.7347408324
$ws.Cells.Item(31, 2).Value = 113825.239530532
$ws.Cells.Item(31, 3).Value = 94931.4685967098
$ws.Cells.Item(31, 4).Value = 86558.4894841983
$ws.Cells.Item(31, 5).Value = 132781.01975517
$ws.Cells.Item(31, 6).Value = 145651.313094884
$ws.Cells.Item(31, 9).Value = 6979.23953053202
$ws.Cells.Item(32, 2).Value = 118293.175213569
$ws.Cells.Item(32, 3).Value = 100358.701554955
$ws.Cells.Item(32, 4).Value = 90442.2662450934
$ws.Cells.Item(32, 5).Value = 138272.584222967
$ws.Cells.Item(32, 6).Value = 150449.447654483
$ws.Cells.Item(32, 9).Value = 5645.17521356934
$ws.Cells.Item(33, 2).Value = 117363.919657343
$ws.Cells.Item(33, 3).Value = 98809.2527911697
$ws.Cells.Item(33, 4).Value = 86149.5061026662
$ws.Cells.Item(33, 5).Value = 136679.161806514
$ws.Cells.Item(33, 6).Value = 147900.562807168
$ws.Cells.Item(33, 9).Value = 1988.91965734278
$ws.Cells.Item(34, 2).Value = 107711.410484044
$ws.Cells.Item(34, 3).Value = 87701.6948230203
$ws.Cells.Item(34, 4).Value = 80219.5866247444
$ws.Cells.Item(34, 5).Value = 127621.608079436
$ws.Cells.Item(34, 6).Value = 138696.040424187
$ws.Cells.Item(34, 9).Value = 6628.41048404369
$ws.Cells.Item(35, 2).Value = 110000.481049071
$ws.Cells.Item(35, 3).Value = 90278.3614508221
$ws.Cells.Item(35, 4).Value = 79850.760307659
$ws.Cells.Item(35, 5).Value = 130255.92780731
$ws.Cells.Item(35, 6).Value = 143319.053169368
$ws.Cells.Item(35, 9).Value = 19985.4810490709
$ws.Cells.Item(36, 2).Value = 109475.025866729
$ws.Cells.Item(36, 3).Value = 89144.433930578
$ws.Cells.Item(36, 4).Value = 77202.9751959559
$ws.Cells.Item(36, 5).Value = 129304.192584712
$ws.Cells.Item(36, 6).Value = 144323.074267837
$ws.Cells.Item(36, 9).Value = 23104.0258667294
